# Auto-generated update of 2024 (column K) violent-crime-ytd figures
# Commit: Add data for 2024-06-11

$wb = $excel.ActiveWorkbook

$updates = [ordered]@{
    'Citywide Totals' = @{ 'K2' = 3376; 'K3' = 3353; 'K4' = 703; 'K5' = 221; 'K6' = 3943; 'K7' = 11596 }
    'Logan Square' = @{ 'K6' = 78; 'K7' = 153 }
    'Austin' = @{ 'K2' = 222; 'K3' = 232; 'K4' = 43; 'K6' = 252; 'K7' = 770 }
    'South Chicago' = @{ 'K6' = 60; 'K7' = 249 }
    'Garfield Park' = @{ 'K3' = 173; 'K7' = 462 }
    'West Pullman' = @{ 'K3' = 67; 'K7' = 188 }
    'Grand Crossing' = @{ 'K2' = 108; 'K6' = 124; 'K7' = 403 }
    'New City' = @{ 'K3' = 73; 'K4' = 9; 'K6' = 109; 'K7' = 281 }
    'By Neighborhood' = @{ 'K7' = 329; 'K8' = 770; 'K9' = 48; 'K15' = 118; 'K19' = 356; 'K20' = 266; 'K21' = 31; 'K23' = 112; 'K25' = 48; 'K29' = 606; 'K33' = 462; 'K37' = 403; 'K41' = 100; 'K42' = 415; 'K43' = 104; 'K44' = 108; 'K46' = 23; 'K47' = 63; 'K51' = 137; 'K52' = 316; 'K53' = 153; 'K54' = 229; 'K58' = 4; 'K63' = 40; 'K65' = 281; 'K67' = 452; 'K72' = 55; 'K73' = 103; 'K78' = 146; 'K79' = 299; 'K83' = 249; 'K84' = 82; 'K85' = 543; 'K86' = 77; 'K88' = 138; 'K89' = 154; 'K90' = 102; 'K94' = 143; 'K95' = 188; 'K96' = 142; 'K97' = 101; 'K101' = 11596 }
    'North Lawndale' = @{ 'K2' = 140; 'K3' = 147; 'K6' = 132; 'K7' = 452 }
    'South Deering' = @{ 'K3' = 31; 'K7' = 82 }
    'Loop' = @{ 'K2' = 41; 'K3' = 68; 'K6' = 107; 'K7' = 229 }
    'Englewood' = @{ 'K3' = 206; 'K6' = 185; 'K7' = 606 }
    'Chatham' = @{ 'K2' = 121; 'K3' = 95; 'K7' = 356 }
    'Irving Park' = @{ 'K2' = 20; 'K6' = 51; 'K7' = 108 }
    'Hermosa' = @{ 'K6' = 42; 'K7' = 100 }
    'Humboldt Park' = @{ 'K2' = 110; 'K3' = 133; 'K7' = 415 }
    'Rogers Park' = @{ 'K2' = 44; 'K7' = 146 }
    'Jefferson Park' = @{ 'K2' = 9; 'K7' = 23 }
    'Douglas' = @{ 'K6' = 29; 'K7' = 112 }
    'West Ridge' = @{ 'K3' = 23; 'K7' = 142 }
    'Chinatown' = @{ 'K6' = 17; 'K7' = 31 }
    'Roseland' = @{ 'K3' = 103; 'K7' = 299 }
    'Near South Side' = @{ 'K3' = 24; 'K6' = 27 }
    'Chicago Lawn' = @{ 'K2' = 93; 'K3' = 75; 'K7' = 266 }
    'Auburn Gresham' = @{ 'K2' = 120; 'K3' = 99; 'K7' = 329 }
    'West Loop' = @{ 'K2' = 40; 'K7' = 143 }
    'East Side' = @{ 'K6' = 8; 'K7' = 48 }
    'Kenwood' = @{ 'K2' = 19; 'K7' = 63 }
    'Brighton Park' = @{ 'K2' = 42; 'K7' = 118 }
    'Avalon Park' = @{ 'K6' = 11; 'K7' = 48 }
    'Portage Park' = @{ 'K6' = 42; 'K7' = 103 }
    'West Town' = @{ 'K6' = 62; 'K7' = 101 }
    'United Center' = @{ 'K3' = 37; 'K7' = 138 }
    'Uptown' = @{ 'K3' = 50; 'K7' = 154 }
    'Streeterville' = @{ 'K4' = 29; 'K7' = 77 }
    'Washington Heights' = @{ 'K4' = 9; 'K7' = 102 }
    'Little Italy, UIC' = @{ 'K5' = 3; 'K7' = 137 }
    'Hyde Park' = @{ 'K2' = 21; 'K6' = 44; 'K7' = 104 }
    'South Shore' = @{ 'K6' = 126; 'K7' = 543 }
    'Old Town' = @{ 'K4' = 4; 'K7' = 55 }
    'Little Village' = @{ 'K6' = 127; 'K7' = 316 }
    'Millenium Park' = @{ 'K6' = 4; 'K7' = 4 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
